$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 405.84616
$ws.Range("I33").Value = 418.9091
$ws.Range("K33").Value = 418.9091
$ws.Range("M33").Value = -189.9091
$ws.Range("H107").Value = 55557988
$ws.Range("I107").Value = 62501236
$ws.Range("K107").Value = 62501236
$ws.Range("M107").Value = -62499316
$ws.Range("H135").Value = 22729928
$ws.Range("I135").Value = 35716490
$ws.Range("K135").Value = 321448410
$ws.Range("M135").Value = -321445875
$ws.Range("H138").Value = 15389150
$ws.Range("I138").Value = 37038610
$ws.Range("J138").Value = 6641.5527
$ws.Range("K138").Value = 111115830
$ws.Range("L138").Value = 19924.6581
$ws.Range("M138").Value = -111110690
$ws.Range("N138").Value = -30204.6581

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3531.1365
$ws.Range("I32").Value = 2734.6
$ws.Range("J32").Value = 11496.5
$ws.Range("K32").Value = 2734.6
$ws.Range("L32").Value = 11496.5
$ws.Range("M32").Value = -2447.6
$ws.Range("N32").Value = -12070.5
$ws.Range("H97").Value = 2575.5151
$ws.Range("I97").Value = 2957
$ws.Range("K97").Value = 2957
$ws.Range("M97").Value = -2461

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H48").Value = 99999
$ws.Range("J48").Value = 99999
$ws.Range("L48").Value = 99999
$ws.Range("N48").Value = -100829
$ws.Range("H134").Value = 4883.1313
$ws.Range("I134").Value = 3751.7354
$ws.Range("K134").Value = 11255.2062
$ws.Range("M134").Value = -8720.206200000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4377.643
$ws.Range("I16").Value = 3467.4285
$ws.Range("J16").Value = 5287.857
$ws.Range("K16").Value = 3467.4285
$ws.Range("L16").Value = 5287.857
$ws.Range("M16").Value = -3180.4285
$ws.Range("N16").Value = -5861.857
$ws.Range("H99").Value = 14327.667
$ws.Range("I99").Value = 5999
$ws.Range("J99").Value = 17103.889
$ws.Range("K99").Value = 5999
$ws.Range("L99").Value = 17103.889
$ws.Range("M99").Value = -4501
$ws.Range("N99").Value = -20099.889
$ws.Range("H113").Value = 4377.643
$ws.Range("I113").Value = 3467.4285
$ws.Range("J113").Value = 5287.857
$ws.Range("K113").Value = 3467.4285
$ws.Range("L113").Value = 5287.857
$ws.Range("M113").Value = -1297.4285
$ws.Range("N113").Value = -9627.857
$ws.Range("H122").Value = 6267.5
$ws.Range("I122").Value = 2173
$ws.Range("J122").Value = 8724.200000000001
$ws.Range("K122").Value = 6519
$ws.Range("L122").Value = 26172.6
$ws.Range("M122").Value = -4069
$ws.Range("N122").Value = -31072.6
$ws.Range("H126").Value = 14327.667
$ws.Range("I126").Value = 5999
$ws.Range("J126").Value = 17103.889
$ws.Range("K126").Value = 17997
$ws.Range("L126").Value = 51311.667
$ws.Range("M126").Value = -15527
$ws.Range("N126").Value = -56251.667
$ws.Range("H132").Value = 10249.875
$ws.Range("I132").Value = 14655.667
$ws.Range("J132").Value = 4585.2856
$ws.Range("K132").Value = 43967.001
$ws.Range("L132").Value = 13755.8568
$ws.Range("M132").Value = -41437.001
$ws.Range("N132").Value = -18815.8568
$ws.Range("H134").Value = 8771.522999999999
$ws.Range("I134").Value = 8937.091
$ws.Range("J134").Value = 8589.4
$ws.Range("K134").Value = 26811.273
$ws.Range("L134").Value = 25768.2
$ws.Range("M134").Value = -24276.273
$ws.Range("N134").Value = -30838.2
$ws.Range("H138").Value = 88667
$ws.Range("J138").Value = 88667
$ws.Range("L138").Value = 88667
$ws.Range("N138").Value = -98947

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 91.875
$ws.Range("J12").Value = 104.85714
$ws.Range("L12").Value = 314.57142
$ws.Range("N12").Value = -660.57142
$ws.Range("H17").Value = 822.375
$ws.Range("I17").Value = 1521.625
$ws.Range("J17").Value = 123.125
$ws.Range("K17").Value = 4564.875
$ws.Range("L17").Value = 369.375
$ws.Range("M17").Value = -4395.875
$ws.Range("N17").Value = -707.375
$ws.Range("H34").Value = 1027.4
$ws.Range("I34").Value = 737.8333
$ws.Range("K34").Value = 2213.4999
$ws.Range("M34").Value = -2129.4999
$ws.Range("H39").Value = 7027.6
$ws.Range("J39").Value = 7971.2856
$ws.Range("L39").Value = 23913.8568
$ws.Range("N39").Value = -24501.8568
$ws.Range("H55").Value = 1766.8889
$ws.Range("J55").Value = 1499
$ws.Range("L55").Value = 4497
$ws.Range("N55").Value = -4851
$ws.Range("H132").Value = 1556.375
$ws.Range("I132").Value = 1487.8182
$ws.Range("J132").Value = 1707.2
$ws.Range("K132").Value = 13390.3638
$ws.Range("L132").Value = 15364.8
$ws.Range("M132").Value = -10860.3638
$ws.Range("N132").Value = -20424.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7423331.5
$ws.Range("I80").Value = 25557966
$ws.Range("K80").Value = 25557966
$ws.Range("M80").Value = -25556968
$ws.Range("H83").Value = 7423331.5
$ws.Range("I83").Value = 25557966
$ws.Range("K83").Value = 127789830
$ws.Range("M83").Value = -127784838
$ws.Range("H140").Value = 137984.28
$ws.Range("J140").Value = 137984.28
$ws.Range("L140").Value = 137984.28
$ws.Range("N140").Value = -148344.28

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4083.36
$ws.Range("I40").Value = 3241.6316
$ws.Range("J40").Value = 6748.8335
$ws.Range("K40").Value = 3241.6316
$ws.Range("L40").Value = 6748.8335
$ws.Range("M40").Value = -3105.6316
$ws.Range("N40").Value = -7020.8335
$ws.Range("H46").Value = 2985.4
$ws.Range("I46").Value = 1231.3334
$ws.Range("J46").Value = 4154.778
$ws.Range("K46").Value = 1231.3334
$ws.Range("L46").Value = 4154.778
$ws.Range("M46").Value = -1043.3334
$ws.Range("N46").Value = -4530.778
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H136").Value = 3534481.2
$ws.Range("I136").Value = 6209536.5
$ws.Range("J136").Value = 8272.137000000001
$ws.Range("K136").Value = 18628609.5
$ws.Range("L136").Value = 24816.411
$ws.Range("M136").Value = -18626059.5
$ws.Range("N136").Value = -29916.411

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 11326.85
$ws.Range("I122").Value = 5039.2144
$ws.Range("J122").Value = 25998
$ws.Range("K122").Value = 15117.6432
$ws.Range("L122").Value = 77994
$ws.Range("M122").Value = -12667.6432
$ws.Range("N122").Value = -82894
$ws.Range("H132").Value = 3143.6135
$ws.Range("I132").Value = 3143.6135
$ws.Range("K132").Value = 9430.8405
$ws.Range("M132").Value = -6900.8405
$ws.Range("H139").Value = 69948
$ws.Range("J139").Value = 69948
$ws.Range("L139").Value = 69948
$ws.Range("N139").Value = -80228
